$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date bump
$ws.Range("A1").Value = 45309

# Price updates (column D)
$updates = @{
    "D19" = 2087.754
    "D20" = 2287.052
    "D21" = 2368.064
    "D22" = 2467.352
    "D23" = 2496.392
    "D24" = 3255.734
    "D25" = 3668.172
    "D26" = 4372.601
    "D27" = 4915.716
    "D28" = 5350.785
    "D29" = 5847.168
    "D30" = 6482.267
    "D31" = 7007.856
    "D32" = 7343.64
    "D33" = 8343.731
    "D35" = 9183.198
    "D36" = 10292.774
    "D37" = 11022.761
    "D43" = 3333.185
    "D44" = 3584.225
    "D45" = 3644.087
    "D46" = 3868.911
    "D47" = 4241.202
    "D48" = 4744.896
    "D49" = 5679.276
    "D50" = 6606.352
    "D51" = 7343.64
    "D52" = 8095.53
    "D53" = 9015.303
    "D54" = 10431.471
    "D55" = 11095.76
    "D56" = 12519.225
    "D57" = 14008.392
    "D58" = 16015.851
    "D59" = 17526.921
    "D68" = 5350.785
    "D69" = 5679.276
    "D70" = 6095.367
    "D71" = 7314.439
    "D72" = 8015.227
    "D73" = 10095.681
    "D74" = 10613.958
    "D75" = 12103.139
    "D76" = 13088.622
    "D77" = 14760.281
    "D78" = 15920.952
    "D79" = 17526.921
    "D80" = 18519.695
    "D81" = 20877.551
    "D82" = 22855.79
    "D83" = 24527.469
    "D90" = 16193.52
    "D91" = 37156.187
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
